$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/border/centered) to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J12
$values = @(
    @(5, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(4, 4),
    @(4, 4),
    @(4, 4),
    @(4, 4),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
